$d = $word.ActiveDocument

# --- Part 1: modify the first paragraph ---
# Original: "This is a Microsoft word document."
# Target:
#   Run 1 (default formatting): "This is a Microsoft word document.  " (two trailing spaces)
#   Run 2 (red C00000):         "(This is a change - Ve"
#   Run 3 (red C00000):         "rsion for branch alternate"
#   Run 4 (red C00000):         ")"
$p1 = $d.Paragraphs(1)
$r1 = $p1.Range
$insertPos = $r1.End - 1   # position right before the paragraph mark

# Append the two trailing spaces to the existing (black) run
$rSpaces = $d.Range($insertPos, $insertPos)
[void]$rSpaces.InsertAfter("  ")
$insertPos = $insertPos + 2

# Red run 1
$seg1 = "(This is a change " + [char]0x2013 + " Ve"
$rSeg1 = $d.Range($insertPos, $insertPos)
[void]$rSeg1.InsertAfter($seg1)
$rColor1 = $d.Range($insertPos, $insertPos + $seg1.Length)
$rColor1.Font.Color = 192
$insertPos = $insertPos + $seg1.Length

# Red run 2
$seg2 = "rsion for branch alternate"
$rSeg2 = $d.Range($insertPos, $insertPos)
[void]$rSeg2.InsertAfter($seg2)
$rColor2 = $d.Range($insertPos, $insertPos + $seg2.Length)
$rColor2.Font.Color = 192
$insertPos = $insertPos + $seg2.Length

# Red run 3
$seg3 = ")"
$rSeg3 = $d.Range($insertPos, $insertPos)
[void]$rSeg3.InsertAfter($seg3)
$rColor3 = $d.Range($insertPos, $insertPos + $seg3.Length)
$rColor3.Font.Color = 192
$insertPos = $insertPos + $seg3.Length

# --- Part 2: append a new empty shaded paragraph at the very end of the body ---
$endRange = $d.Content
$endRange.Collapse(0)
$xmlFrag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$endRange.InsertXML($xmlFrag)

Write-Output "done"
